# Update TPM-derived NATMI ligand-receptor statistics (C3-Cd46)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.820425
$ws.Range("H2").Value = 11.461275
$ws.Range("I2").Value = 0.02049663039797357
$ws.Range("J2").Value = 0.02049663039797357
$ws.Range("M2").Value = 3.477616333333334
$ws.Range("N2").Value = 10.432849
$ws.Range("O2").Value = 0.4986219472823802
$ws.Range("P2").Value = 0.4986219472823803
$ws.Range("Q2").Value = 13.285972380275
$ws.Range("R2").Value = 119.573751422475
$ws.Range("S2").Value = 0.01022006976176481
$ws.Range("T2").Value = 0.01022006976176481
# Row 3
$ws.Range("G3").Value = 3.820425
$ws.Range("H3").Value = 11.461275
$ws.Range("I3").Value = 0.02049663039797357
$ws.Range("J3").Value = 0.02049663039797357
$ws.Range("O3").Value = 0.2262779433422874
$ws.Range("P3").Value = 0.2262779433422874
$ws.Range("Q3").Value = 6.0292622936
$ws.Range("R3").Value = 54.2633606424
$ws.Range("S3").Value = 0.004637935371900469
$ws.Range("T3").Value = 0.004637935371900469
# Row 4
$ws.Range("G4").Value = 3.820425
$ws.Range("H4").Value = 11.461275
$ws.Range("I4").Value = 0.02049663039797357
$ws.Range("J4").Value = 0.02049663039797357
$ws.Range("M4").Value = 0.3701243333333333
$ws.Range("N4").Value = 1.110373
$ws.Range("O4").Value = 0.05306856712579453
$ws.Range("P4").Value = 0.05306856712579454
$ws.Range("Q4").Value = 1.414032256175
$ws.Range("R4").Value = 12.726290305575
$ws.Range("S4").Value = 0.001087726806127461
$ws.Range("T4").Value = 0.001087726806127461
# Row 5
$ws.Range("G5").Value = 3.820425
$ws.Range("H5").Value = 11.461275
$ws.Range("I5").Value = 0.02049663039797357
$ws.Range("J5").Value = 0.02049663039797357
$ws.Range("M5").Value = 0.775099
$ws.Range("N5").Value = 2.325297
$ws.Range("O5").Value = 0.1111339882471103
$ws.Range("P5").Value = 0.1111339882471103
$ws.Range("Q5").Value = 2.961207597075
$ws.Range("R5").Value = 26.650868373675
$ws.Range("S5").Value = 0.002277872281753759
$ws.Range("T5").Value = 0.002277872281753759
# Row 6
$ws.Range("G6").Value = 3.820425
$ws.Range("H6").Value = 11.461275
$ws.Range("I6").Value = 0.02049663039797357
$ws.Range("J6").Value = 0.02049663039797357
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.7734500000000001
$ws.Range("N6").Value = 2.32035
$ws.Range("O6").Value = 0.1108975540024274
$ws.Range("P6").Value = 0.1108975540024274
$ws.Range("Q6").Value = 2.954907716250001
$ws.Range("R6").Value = 26.59416944625001
$ws.Range("S6").Value = 0.002273026176427069
$ws.Range("T6").Value = 0.002273026176427069
# Row 7
$ws.Range("I7").Value = 0.7542622677884155
$ws.Range("J7").Value = 0.7542622677884157
$ws.Range("M7").Value = 3.477616333333334
$ws.Range("N7").Value = 10.432849
$ws.Range("O7").Value = 0.4986219472823802
$ws.Range("P7").Value = 0.4986219472823803
$ws.Range("Q7").Value = 488.9148832147176
$ws.Range("R7").Value = 4400.233948932459
$ws.Range("S7").Value = 0.3760917207262839
$ws.Range("T7").Value = 0.376091720726284
# Row 8
$ws.Range("I8").Value = 0.7542622677884155
$ws.Range("J8").Value = 0.7542622677884157
$ws.Range("O8").Value = 0.2262779433422874
$ws.Range("P8").Value = 0.2262779433422874
$ws.Range("S8").Value = 0.1706729146958523
$ws.Range("T8").Value = 0.1706729146958524
# Row 9
$ws.Range("I9").Value = 0.7542622677884155
$ws.Range("J9").Value = 0.7542622677884157
$ws.Range("M9").Value = 0.3701243333333333
$ws.Range("N9").Value = 1.110373
$ws.Range("O9").Value = 0.05306856712579453
$ws.Range("P9").Value = 0.05306856712579454
$ws.Range("Q9").Value = 52.03543975569622
$ws.Range("R9").Value = 468.318957801266
$ws.Range("S9").Value = 0.04002761778858354
$ws.Range("T9").Value = 0.04002761778858355
# Row 10
$ws.Range("I10").Value = 0.7542622677884155
$ws.Range("J10").Value = 0.7542622677884157
$ws.Range("M10").Value = 0.775099
$ws.Range("N10").Value = 2.325297
$ws.Range("O10").Value = 0.1111339882471103
$ws.Range("P10").Value = 0.1111339882471103
$ws.Range("Q10").Value = 108.9704558356527
$ws.Range("R10").Value = 980.734102520874
$ws.Range("S10").Value = 0.08382417400363656
$ws.Range("T10").Value = 0.08382417400363658
# Row 11
$ws.Range("I11").Value = 0.7542622677884155
$ws.Range("J11").Value = 0.7542622677884157
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.7734500000000001
$ws.Range("N11").Value = 2.32035
$ws.Range("O11").Value = 0.1108975540024274
$ws.Range("P11").Value = 0.1108975540024274
$ws.Range("Q11").Value = 108.7386244416333
$ws.Range("R11").Value = 978.6476199747002
$ws.Range("S11").Value = 0.08364584057405919
$ws.Range("T11").Value = 0.08364584057405922
# Row 12
$ws.Range("G12").Value = 30.51067
$ws.Range("H12").Value = 91.53201
$ws.Range("I12").Value = 0.1636901460399144
$ws.Range("J12").Value = 0.1636901460399144
$ws.Range("M12").Value = 3.477616333333334
$ws.Range("N12").Value = 10.432849
$ws.Range("O12").Value = 0.4986219472823802
$ws.Range("P12").Value = 0.4986219472823803
$ws.Range("Q12").Value = 106.1044043329433
$ws.Range("R12").Value = 954.9396389964901
$ws.Range("S12").Value = 0.08161949936935935
$ws.Range("T12").Value = 0.08161949936935935
# Row 13
$ws.Range("G13").Value = 30.51067
$ws.Range("H13").Value = 91.53201
$ws.Range("I13").Value = 0.1636901460399144
$ws.Range("J13").Value = 0.1636901460399144
$ws.Range("O13").Value = 0.2262779433422874
$ws.Range("P13").Value = 0.2262779433422874
$ws.Range("Q13").Value = 48.15088169077334
$ws.Range("R13").Value = 433.35793521696
$ws.Range("S13").Value = 0.03703946959131051
$ws.Range("T13").Value = 0.03703946959131051
# Row 14
$ws.Range("G14").Value = 30.51067
$ws.Range("H14").Value = 91.53201
$ws.Range("I14").Value = 0.1636901460399144
$ws.Range("J14").Value = 0.1636901460399144
$ws.Range("M14").Value = 0.3701243333333333
$ws.Range("N14").Value = 1.110373
$ws.Range("O14").Value = 0.05306856712579453
$ws.Range("P14").Value = 0.05306856712579454
$ws.Range("Q14").Value = 11.29274139330333
$ws.Range("R14").Value = 101.63467253973
$ws.Range("S14").Value = 0.008686801502950309
$ws.Range("T14").Value = 0.008686801502950311
# Row 15
$ws.Range("G15").Value = 30.51067
$ws.Range("H15").Value = 91.53201
$ws.Range("I15").Value = 0.1636901460399144
$ws.Range("J15").Value = 0.1636901460399144
$ws.Range("M15").Value = 0.775099
$ws.Range("N15").Value = 2.325297
$ws.Range("O15").Value = 0.1111339882471103
$ws.Range("P15").Value = 0.1111339882471103
$ws.Range("Q15").Value = 23.64878980633
$ws.Range("R15").Value = 212.83910825697
$ws.Range("S15").Value = 0.01819153876616763
$ws.Range("T15").Value = 0.01819153876616763
# Row 16
$ws.Range("G16").Value = 30.51067
$ws.Range("H16").Value = 91.53201
$ws.Range("I16").Value = 0.1636901460399144
$ws.Range("J16").Value = 0.1636901460399144
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.7734500000000001
$ws.Range("N16").Value = 2.32035
$ws.Range("O16").Value = 0.1108975540024274
$ws.Range("P16").Value = 0.1108975540024274
$ws.Range("Q16").Value = 23.5984777115
$ws.Range("R16").Value = 212.3862994035
$ws.Range("S16").Value = 0.01815283681012664
$ws.Range("T16").Value = 0.01815283681012664
# Row 17
$ws.Range("G17").Value = 0.258813
$ws.Range("H17").Value = 0.776439
$ws.Range("I17").Value = 0.001388535150720334
$ws.Range("J17").Value = 0.001388535150720334
$ws.Range("M17").Value = 3.477616333333334
$ws.Range("N17").Value = 10.432849
$ws.Range("O17").Value = 0.4986219472823802
$ws.Range("P17").Value = 0.4986219472823803
$ws.Range("Q17").Value = 0.9000523160790002
$ws.Range("R17").Value = 8.100470844711001
$ws.Range("S17").Value = 0.0006923541007222064
$ws.Range("T17").Value = 0.0006923541007222064
# Row 18
$ws.Range("G18").Value = 0.258813
$ws.Range("H18").Value = 0.776439
$ws.Range("I18").Value = 0.001388535150720334
$ws.Range("J18").Value = 0.001388535150720334
$ws.Range("O18").Value = 0.2262779433422874
$ws.Range("P18").Value = 0.2262779433422874
$ws.Range("Q18").Value = 0.408449704416
$ws.Range("R18").Value = 3.676047339744
$ws.Range("S18").Value = 0.0003141948781634703
$ws.Range("T18").Value = 0.0003141948781634703
# Row 19
$ws.Range("G19").Value = 0.258813
$ws.Range("H19").Value = 0.776439
$ws.Range("I19").Value = 0.001388535150720334
$ws.Range("J19").Value = 0.001388535150720334
$ws.Range("M19").Value = 0.3701243333333333
$ws.Range("N19").Value = 1.110373
$ws.Range("O19").Value = 0.05306856712579453
$ws.Range("P19").Value = 0.05306856712579454
$ws.Range("Q19").Value = 0.095792989083
$ws.Range("R19").Value = 0.8621369017470001
$ws.Range("S19").Value = 0.00007368757085252728
$ws.Range("T19").Value = 0.00007368757085252728
# Row 20
$ws.Range("G20").Value = 0.258813
$ws.Range("H20").Value = 0.776439
$ws.Range("I20").Value = 0.001388535150720334
$ws.Range("J20").Value = 0.001388535150720334
$ws.Range("M20").Value = 0.775099
$ws.Range("N20").Value = 2.325297
$ws.Range("O20").Value = 0.1111339882471103
$ws.Range("P20").Value = 0.1111339882471103
$ws.Range("Q20").Value = 0.200605697487
$ws.Range("R20").Value = 1.805451277383
$ws.Range("S20").Value = 0.0001543134491208532
$ws.Range("T20").Value = 0.0001543134491208532
# Row 21
$ws.Range("G21").Value = 0.258813
$ws.Range("H21").Value = 0.776439
$ws.Range("I21").Value = 0.001388535150720334
$ws.Range("J21").Value = 0.001388535150720334
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 0.7734500000000001
$ws.Range("N21").Value = 2.32035
$ws.Range("O21").Value = 0.1108975540024274
$ws.Range("P21").Value = 0.1108975540024274
$ws.Range("Q21").Value = 0.20017891485
$ws.Range("R21").Value = 1.80161023365
$ws.Range("S21").Value = 0.000153985151861277
$ws.Range("T21").Value = 0.000153985151861277
# Row 22
$ws.Range("G22").Value = 11.213844
$ws.Range("H22").Value = 33.641532
$ws.Range("I22").Value = 0.0601624206229761
$ws.Range("J22").Value = 0.0601624206229761
$ws.Range("M22").Value = 3.477616333333334
$ws.Range("N22").Value = 10.432849
$ws.Range("O22").Value = 0.4986219472823802
$ws.Range("P22").Value = 0.4986219472823803
$ws.Range("Q22").Value = 38.99744705385201
$ws.Range("R22").Value = 350.977023484668
$ws.Range("S22").Value = 0.02999830332424997
$ws.Range("T22").Value = 0.02999830332424997
# Row 23
$ws.Range("G23").Value = 11.213844
$ws.Range("H23").Value = 33.641532
$ws.Range("I23").Value = 0.0601624206229761
$ws.Range("J23").Value = 0.0601624206229761
$ws.Range("O23").Value = 0.2262779433422874
$ws.Range("P23").Value = 0.2262779433422874
$ws.Range("Q23").Value = 17.697299854208
$ws.Range("R23").Value = 159.275698687872
$ws.Range("S23").Value = 0.01361342880506065
$ws.Range("T23").Value = 0.01361342880506065
# Row 24
$ws.Range("G24").Value = 11.213844
$ws.Range("H24").Value = 33.641532
$ws.Range("I24").Value = 0.0601624206229761
$ws.Range("J24").Value = 0.0601624206229761
$ws.Range("M24").Value = 0.3701243333333333
$ws.Range("N24").Value = 1.110373
$ws.Range("O24").Value = 0.05306856712579453
$ws.Range("P24").Value = 0.05306856712579454
$ws.Range("Q24").Value = 4.150516534604
$ws.Range("R24").Value = 37.354648811436
$ws.Range("S24").Value = 0.003192733457280693
$ws.Range("T24").Value = 0.003192733457280693
# Row 25
$ws.Range("G25").Value = 11.213844
$ws.Range("H25").Value = 33.641532
$ws.Range("I25").Value = 0.0601624206229761
$ws.Range("J25").Value = 0.0601624206229761
$ws.Range("M25").Value = 0.775099
$ws.Range("N25").Value = 2.325297
$ws.Range("O25").Value = 0.1111339882471103
$ws.Range("P25").Value = 0.1111339882471103
$ws.Range("Q25").Value = 8.691839270556001
$ws.Range("R25").Value = 78.22655343500399
$ws.Range("S25").Value = 0.006686089746431535
$ws.Range("T25").Value = 0.006686089746431535
# Row 26
$ws.Range("G26").Value = 11.213844
$ws.Range("H26").Value = 33.641532
$ws.Range("I26").Value = 0.0601624206229761
$ws.Range("J26").Value = 0.0601624206229761
$ws.Range("K26").Value = 3
$ws.Range("L26").Value = 1
$ws.Range("M26").Value = 0.7734500000000001
$ws.Range("N26").Value = 2.32035
$ws.Range("O26").Value = 0.1108975540024274
$ws.Range("P26").Value = 0.1108975540024274
$ws.Range("Q26").Value = 8.673347641800001
$ws.Range("R26").Value = 78.0601287762
$ws.Range("S26").Value = 0.006671865289953246
$ws.Range("T26").Value = 0.006671865289953246
